# Update the StructureDefinition-employee-department workbook:
#   - rebrand from "ibm.com / Alvearie" to "linuxforhealth.org / LinuxForHealth"
#   - bump Version 7.0.0 -> 8.0.0
#   - bump Date to the new publish timestamp
#   - clear the stray Constraint(s) text that had leaked onto the root
#     "Extension" row on the Elements sheet (it belongs only on Extension.extension)

$wb = $excel.ActiveWorkbook

$newUrl = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-department"

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = $newUrl
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Extension.url's "Fixed Value" column also held the old URL (it shared
# the same string as Metadata!B2) - keep it in sync with the new URL.
$elements.Range("Q5").Value = $newUrl

# The root Extension row's "Constraint(s)" cell incorrectly duplicated the
# ele-1/ext-1 constraint text that belongs to the Extension.extension row
# (AI4) only; clear it out on the root row.
$elements.Range("AI2").Value = ""
